# Auto-generated edit script applying meteocat daily-summary refresh
# (timestamps + updated readings) for commit "Update automatic: dades i banners [2026-02-18 18:20]"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding a bare "NN%" reading need their number format pinned to
# Text first, otherwise the COM layer helpfully (re)parses the literal
# "82%" string into the numeric percentage 0.82 - exactly like typing it
# into Excel by hand would. The source feed stores these as plain text,
# so we preserve that by forcing @ (Text) before writing the value.
$percentCells = @("H3", "H9", "H10", "H12", "H15", "H16", "H17", "H19", "H20", "H23", "H24", "H27", "H30", "H31", "H32", "H36", "H42", "H43", "H44", "H46")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-18 18:18:52"
$ws.Range("O2").Value = "2.0 °C"
$ws.Range("E3").Value = "2026-02-18 18:18:54"
$ws.Range("H3").Value = "70%"
$ws.Range("O3").Value = "-0.2 °C"
$ws.Range("E4").Value = "2026-02-18 18:18:57"
$ws.Range("J4").Value = "1014.1 hPa"
$ws.Range("L4").Value = "35.3 km/h - 277º 17:52 TU"
$ws.Range("O4").Value = "11.6 °C"
$ws.Range("E5").Value = "2026-02-18 18:19:00"
$ws.Range("E6").Value = "2026-02-18 18:19:03"
$ws.Range("J6").Value = "1013.8 hPa"
$ws.Range("O6").Value = "11.6 °C"
$ws.Range("E7").Value = "2026-02-18 18:19:05"
$ws.Range("J7").Value = "1015.1 hPa"
$ws.Range("E8").Value = "2026-02-18 18:19:08"
$ws.Range("J8").Value = "1014.8 hPa"
$ws.Range("E9").Value = "2026-02-18 18:19:11"
$ws.Range("H9").Value = "82%"
$ws.Range("O9").Value = "10.3 °C"
$ws.Range("E10").Value = "2026-02-18 18:19:14"
$ws.Range("H10").Value = "83%"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-18 18:19:16"
$ws.Range("O11").Value = "5.6 °C"
$ws.Range("E12").Value = "2026-02-18 18:19:19"
$ws.Range("H12").Value = "90%"
$ws.Range("O12").Value = "10.2 °C"
$ws.Range("E13").Value = "2026-02-18 18:19:22"
$ws.Range("J13").Value = "1016.5 hPa"
$ws.Range("O13").Value = "3.5 °C"
$ws.Range("E14").Value = "2026-02-18 18:19:24"
$ws.Range("E15").Value = "2026-02-18 18:19:27"
$ws.Range("H15").Value = "82%"
$ws.Range("M15").Value = "17.1 °C 17:41 TU"
$ws.Range("O15").Value = "10.2 °C"
$ws.Range("E16").Value = "2026-02-18 18:19:30"
$ws.Range("H16").Value = "47%"
$ws.Range("N16").Value = "-2.6 °C 17:43 TU"
$ws.Range("O16").Value = "0.2 °C"
$ws.Range("E17").Value = "2026-02-18 18:19:32"
$ws.Range("H17").Value = "85%"
$ws.Range("E18").Value = "2026-02-18 18:19:35"
$ws.Range("J18").Value = "1014.2 hPa"
$ws.Range("L18").Value = "40.3 km/h - 209º 17:48 TU"
$ws.Range("O18").Value = "11.4 °C"
$ws.Range("E19").Value = "2026-02-18 18:19:38"
$ws.Range("H19").Value = "90%"
$ws.Range("L19").Value = "51.5 km/h - 250º 17:34 TU"
$ws.Range("O19").Value = "7.1 °C"
$ws.Range("E20").Value = "2026-02-18 18:19:40"
$ws.Range("H20").Value = "71%"
$ws.Range("L20").Value = "52.2 km/h - 267º 17:34 TU"
$ws.Range("N20").Value = "-2.5 °C 17:45 TU"
$ws.Range("E21").Value = "2026-02-18 18:19:43"
$ws.Range("J21").Value = "1015.7 hPa"
$ws.Range("O21").Value = "6.3 °C"
$ws.Range("E22").Value = "2026-02-18 18:19:46"
$ws.Range("N22").Value = "-3.8 °C 17:31 TU"
$ws.Range("E23").Value = "2026-02-18 18:19:48"
$ws.Range("H23").Value = "51%"
$ws.Range("N23").Value = "-2.8 °C 17:52 TU"
$ws.Range("E24").Value = "2026-02-18 18:19:51"
$ws.Range("H24").Value = "84%"
$ws.Range("J24").Value = "1015.6 hPa"
$ws.Range("O24").Value = "9.5 °C"
$ws.Range("E25").Value = "2026-02-18 18:19:53"
$ws.Range("L25").Value = "57.6 km/h - 251º 17:57 TU"
$ws.Range("E26").Value = "2026-02-18 18:19:56"
$ws.Range("J26").Value = "1012.9 hPa"
$ws.Range("E27").Value = "2026-02-18 18:19:59"
$ws.Range("H27").Value = "52%"
$ws.Range("N27").Value = "-1.1 °C 17:59 TU"
$ws.Range("O27").Value = "1.9 °C"
$ws.Range("E28").Value = "2026-02-18 18:20:02"
$ws.Range("J28").Value = "1014.0 hPa"
$ws.Range("O28").Value = "9.4 °C"
$ws.Range("E29").Value = "2026-02-18 18:20:04"
$ws.Range("E30").Value = "2026-02-18 18:20:07"
$ws.Range("H30").Value = "81%"
$ws.Range("J30").Value = "1013.7 hPa"
$ws.Range("O30").Value = "10.3 °C"
$ws.Range("E31").Value = "2026-02-18 18:20:09"
$ws.Range("H31").Value = "74%"
$ws.Range("J31").Value = "1012.4 hPa"
$ws.Range("L31").Value = "98.3 km/h - 209º 17:55 TU"
$ws.Range("M31").Value = "16.8 °C 17:48 TU"
$ws.Range("O31").Value = "12.4 °C"
$ws.Range("E32").Value = "2026-02-18 18:20:12"
$ws.Range("H32").Value = "82%"
$ws.Range("E33").Value = "2026-02-18 18:20:15"
$ws.Range("J33").Value = "1015.0 hPa"
$ws.Range("O33").Value = "4.8 °C"
$ws.Range("E34").Value = "2026-02-18 18:20:18"
$ws.Range("E35").Value = "2026-02-18 18:20:21"
$ws.Range("J35").Value = "1014.9 hPa"
$ws.Range("O35").Value = "9.6 °C"
$ws.Range("E36").Value = "2026-02-18 18:20:23"
$ws.Range("H36").Value = "88%"
$ws.Range("J36").Value = "1014.1 hPa"
$ws.Range("O36").Value = "11.4 °C"
$ws.Range("E37").Value = "2026-02-18 18:20:26"
$ws.Range("J37").Value = "1015.7 hPa"
$ws.Range("O37").Value = "5.7 °C"
$ws.Range("E38").Value = "2026-02-18 18:20:29"
$ws.Range("O38").Value = "12.2 °C"
$ws.Range("E39").Value = "2026-02-18 18:20:32"
$ws.Range("L39").Value = "105.8 km/h - 300º 17:57 TU"
$ws.Range("O39").Value = "1.3 °C"
$ws.Range("E40").Value = "2026-02-18 18:20:35"
$ws.Range("J40").Value = "1016.4 hPa"
$ws.Range("O40").Value = "5.9 °C"
$ws.Range("E41").Value = "2026-02-18 18:20:38"
$ws.Range("J41").Value = "1015.3 hPa"
$ws.Range("E42").Value = "2026-02-18 18:20:41"
$ws.Range("H42").Value = "87%"
$ws.Range("O42").Value = "11.4 °C"
$ws.Range("E43").Value = "2026-02-18 18:20:43"
$ws.Range("H43").Value = "81%"
$ws.Range("E44").Value = "2026-02-18 18:20:46"
$ws.Range("H44").Value = "70%"
$ws.Range("E45").Value = "2026-02-18 18:20:49"
$ws.Range("J45").Value = "1012.7 hPa"
$ws.Range("O45").Value = "7.2 °C"
$ws.Range("E46").Value = "2026-02-18 18:20:52"
$ws.Range("H46").Value = "82%"
$ws.Range("J46").Value = "1015.6 hPa"
$ws.Range("O46").Value = "11.0 °C"
